# Jogos_da_Semana_FlashScore_2024-11-12.xlsx update
# --------------------------------------------------
# The source feed re-pulled this week's fixture list. The old first match
# (Guarani x Amazonas, row 2) dropped off the list entirely, the match that
# used to be row 3 (Botafogo SP x Ceara) is now the only/first fixture, and
# its odds (including a reshuffled "correct score" column order and a few
# re-priced odds) were refreshed. Net effect: the sheet shrinks from two
# data rows to one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old row 2 (Guarani x Amazonas). Excel shifts the old row 3
#    (Botafogo SP x Ceara) up into row 2 automatically, carrying its values
#    (id, date, time, league, teams, and most odds) along with it.
$ws.Rows.Item(2).Delete()

# 2) Header row: the "correct score" odds columns AG:AM were reordered so
#    that Odd_CS_4-4 now leads the block (AG), with the remaining six
#    labels sliding one column to the right.
$ws.Range("AG1").Value = "Odd_CS_4-4"
$ws.Range("AH1").Value = "Odd_CS_0-1"
$ws.Range("AI1").Value = "Odd_CS_0-2"
$ws.Range("AJ1").Value = "Odd_CS_1-2"
$ws.Range("AK1").Value = "Odd_CS_0-3"
$ws.Range("AL1").Value = "Odd_CS_1-3"
$ws.Range("AM1").Value = "Odd_CS_2-3"

# 3) Row 2 (now Botafogo SP x Ceara): a handful of odds were re-priced...
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 4.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("AQ2").Value = 81
$ws.Range("AX2").Value = 12
$ws.Range("AY2").Value = 29
$ws.Range("BA2").Value = 81

# ...and the AG:AM correct-score odds follow the same column reorder as the
# header (old AM value moves to AG, old AG:AL slide right to AH:AM).
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 5.5
$ws.Range("AI2").Value = 8
$ws.Range("AJ2").Value = 9.5
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 21
$ws.Range("AM2").Value = 41

Write-Host "done"
